{"js": "// Replace the date heading and every division problem in the practice table\n// with the values from the new day's worksheet. Each (old, new) pair is a\n// unique piece of text in the document, so a simple find/replace per pair\n// is safe and keeps all existing run formatting (font/size) intact.\nconst replacements = [\n  [\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"],\n  [\"431\u00f77=\", \"804\u00f73=\"],\n  [\"395\u00f74=\", \"232\u00f74=\"],\n  [\"706\u00f74=\", \"473\u00f79=\"],\n  [\"335\u00f76=\", \"482\u00f78=\"],\n  [\"196\u00f73=\", \"400\u00f75=\"],\n  [\"874\u00f74=\", \"225\u00f78=\"],\n  [\"159\u00f79=\", \"974\u00f74=\"],\n  [\"260\u00f74=\", \"612\u00f72=\"],\n  [\"373\u00f72=\", \"562\u00f73=\"],\n  [\"530\u00f72=\", \"900\u00f75=\"],\n  [\"728\u00f76=\", \"273\u00f77=\"],\n  [\"693\u00f73=\", \"485\u00f79=\"],\n  [\"151\u00f73=\", \"956\u00f78=\"],\n  [\"228\u00f77=\", \"581\u00f75=\"],\n  [\"499\u00f72=\", \"403\u00f78=\"],\n  [\"799\u00f72=\", \"995\u00f77=\"],\n  [\"348\u00f72=\", \"968\u00f79=\"],\n  [\"167\u00f73=\", \"921\u00f77=\"],\n  [\"174\u00f76=\", \"457\u00f78=\"],\n  [\"784\u00f78=\", \"820\u00f72=\"],\n  [\"668\u00f73=\", \"420\u00f76=\"],\n  [\"449\u00f73=\", \"214\u00f79=\"],\n  [\"276\u00f79=\", \"965\u00f79=\"],\n  [\"429\u00f72=\", \"679\u00f72=\"],\n  [\"397\u00f75=\", \"499\u00f77=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the date heading and every division problem in the practice table\n# with the values from the new day's worksheet. Each (old, new) pair is a\n# unique piece of text in the document, so a Find/Replace pass per pair is\n# safe and leaves all existing run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"),\n    @(\"431\u00f77=\", \"804\u00f73=\"),\n    @(\"395\u00f74=\", \"232\u00f74=\"),\n    @(\"706\u00f74=\", \"473\u00f79=\"),\n    @(\"335\u00f76=\", \"482\u00f78=\"),\n    @(\"196\u00f73=\", \"400\u00f75=\"),\n    @(\"874\u00f74=\", \"225\u00f78=\"),\n    @(\"159\u00f79=\", \"974\u00f74=\"),\n    @(\"260\u00f74=\", \"612\u00f72=\"),\n    @(\"373\u00f72=\", \"562\u00f73=\"),\n    @(\"530\u00f72=\", \"900\u00f75=\"),\n    @(\"728\u00f76=\", \"273\u00f77=\"),\n    @(\"693\u00f73=\", \"485\u00f79=\"),\n    @(\"151\u00f73=\", \"956\u00f78=\"),\n    @(\"228\u00f77=\", \"581\u00f75=\"),\n    @(\"499\u00f72=\", \"403\u00f78=\"),\n    @(\"799\u00f72=\", \"995\u00f77=\"),\n    @(\"348\u00f72=\", \"968\u00f79=\"),\n    @(\"167\u00f73=\", \"921\u00f77=\"),\n    @(\"174\u00f76=\", \"457\u00f78=\"),\n    @(\"784\u00f78=\", \"820\u00f72=\"),\n    @(\"668\u00f73=\", \"420\u00f76=\"),\n    @(\"449\u00f73=\", \"214\u00f79=\"),\n    @(\"276\u00f79=\", \"965\u00f79=\"),\n    @(\"429\u00f72=\", \"679\u00f72=\"),\n    @(\"397\u00f75=\", \"499\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
